$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ToDo item: "Admin panel(crud system for subsribers)" goes into row 26,
# continuing the existing A:list / B:checkbox-column pattern (last existing
# row was 25 = "newsletter with mailtrap").
$ws.Range("A26").Value = "Admin panel(crud system for subsribers)"

# Copy the formatting (fill/font cell style) from the row above (A25/B25)
# down onto the new row so the new row re-uses the existing style indices
# instead of creating brand-new ones.
$ws.Range("A25:B25").Copy()
$ws.Range("A26:B26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection to the new last row, matching where Excel
# leaves the cursor after typing the new entry.
$ws.Range("A26").Select()
